$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A values (text) in the order needed to match shared string table order
$ws.Cells.Item(3, 1).Value = "bearing_61814_inner"
$ws.Cells.Item(4, 1).Value = "bearing_61814_outer"
$ws.Cells.Item(8, 1).Value = "M3_adjustable"
$ws.Cells.Item(5, 1).Value = "bearing_61705_outer"
$ws.Cells.Item(6, 1).Value = "bearing_61705_inner"
$ws.Cells.Item(7, 1).Value = "bearing_606_outer"
$ws.Cells.Item(10, 1).Value = "gearbox_42_outer"
$ws.Cells.Item(9, 1).Value = "gearbox_35_outer"

# Now set column B values
$ws.Cells.Item(3, 2).Value = 70
$ws.Cells.Item(4, 2).Value = 90.1
$ws.Cells.Item(5, 2).Value = 32.1
$ws.Cells.Item(6, 2).Value = 25
$ws.Cells.Item(7, 2).Value = 17.1
$ws.Cells.Item(8, 2).Value = 3.2
$ws.Cells.Item(9, 2).Value = 32.3
$ws.Cells.Item(10, 2).Value = 42.3

$ws.Range("B11").Select()
